# "The register is now work"
# Populate the "Registro" sheet with its first movement entries, rename the
# ITEM header to EQUIPAMENTO, widen column B to fit it, set the print setup,
# and leave the selection where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registro")

# Header tweak: "ITEM" -> "EQUIPAMENTO"
$ws.Range("B2").Value = "EQUIPAMENTO"

# Widen column B so the new header fits (~19.43 characters once stored)
$ws.Columns.Item(2).ColumnWidth = 18.6

# Row 3: first movement record (a keyboard going out)
$ws.Range("A3").Value = "NULO"
$ws.Range("B3").Value = "TECLADO"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "SAIDA"
$ws.Range("E3").Value = "15/5/2024"

# Row 4: second movement record (a keyboard coming back in)
$ws.Range("A4").Value = "NULO"
$ws.Range("B4").Value = "TECLADO"
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "ENTRADA"
$ws.Range("E4").Value = "15/5/2024"
$ws.Range("E4").NumberFormat = "m/d/yyyy"

# Print setup for the sheet
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# Leave the selection where the user ended up
$ws.Range("E10").Select() | Out-Null
